$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 112, shifting existing rows 112-116 down to 113-117
$ws.Rows.Item(112).Insert()

# Populate the new row 112 with the new weekly data
$ws.Range("A112").Value = 11
$ws.Range("B112").Value = "Vega Monumental Concepción"
$ws.Range("C112").Value = "Bíobío"
$ws.Range("D112").Value = 44747
$ws.Range("E112").Value = 8
$ws.Range("F112").Value = 100112021
$ws.Range("G112").Value = "Ají"
$ws.Range("H112").Value = "Americana (o)"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 22
$ws.Range("K112").Value = 11000
$ws.Range("L112").Value = 12000
$ws.Range("M112").Value = 11545
$ws.Range("N112").Value = "$/caja 15 kilos"
$ws.Range("O112").Value = "Provincia de Huasco"
$ws.Range("P112").Value = 770
$ws.Range("Q112").Value = 15
$ws.Range("R112").Value = "Hortaliza"
